# Add the new HealthNav medical-plan Q&A row to the questions table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "What does my medical plan cover?"
$ws.Range("B4").Value = "Your medical plan covers office visits, specialist visits, outpatient and inpatient services, diagnostic services, emergency care, physical therapy, mental health services, and prescription drugs through BCBS HDHP network."
